# This script inserts 5 new price-report rows (new rows 319-323) into the
# "Tomate" sheet, pushing the previously-existing rows 319-384 down to
# 324-389. The new rows contain a new reporting date (44474) worth of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 319 - this shifts rows 319:384 down
# to 324:389 and keeps all their original values/formatting intact.
$ws.Rows("319:323").Insert()

# Common (unchanged across the 5 new rows) column values.
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$categoriaId = 100112020
$categoria = "Tomate"
$variedad = "Larga vida"
$clasificacion = "Hortaliza"
$fecha = 44474

$newRows = @(
    @{ Row = 319; Calidad = "Extra";   Volumen = 300; PMin = 25000; PMax = 25000; PProm = 25000; Unidad = "$/bandeja 18 kilos"; Origen = "Región de Arica y Parinacota"; PKg = 1389; Kg = 18 },
    @{ Row = 320; Calidad = "Extra";   Volumen = 300; PMin = 28000; PMax = 28000; PProm = 28000; Unidad = "$/bandeja 20 kilos"; Origen = "Región de Arica y Parinacota"; PKg = 1400; Kg = 20 },
    @{ Row = 321; Calidad = "Primera"; Volumen = 300; PMin = 24000; PMax = 24000; PProm = 24000; Unidad = "$/bandeja 18 kilos"; Origen = "Región de Arica y Parinacota"; PKg = 1333; Kg = 18 },
    @{ Row = 322; Calidad = "Primera"; Volumen = 300; PMin = 27000; PMax = 27000; PProm = 27000; Unidad = "$/bandeja 20 kilos"; Origen = "Región de Arica y Parinacota"; PKg = 1350; Kg = 20 },
    @{ Row = 323; Calidad = "Segunda"; Volumen = 300; PMin = 21000; PMax = 21000; PProm = 21000; Unidad = "$/bandeja 18 kilos"; Origen = "Región de Arica y Parinacota"; PKg = 1167; Kg = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $r.Unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $r.Kg
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
